$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: row 9 SmartScore cells were stored as text; fix them up as real numbers ---
$ws.Range("I9").Value = 0.71
$ws.Range("L9").Value = 0.438
$ws.Range("O9").Value = 0.429
$ws.Range("R9").Value = 0.761
$ws.Range("U9").Value = 0.706
$ws.Range("X9").Value = 0.517
$ws.Range("AA9").Value = 0.67
$ws.Range("AD9").Value = 0.582
$ws.Range("AG9").Value = 0.447

# --- Step 2: append new row 10 (Ilse Aguirre submission from Streamlit) ---
$ws.Range("A10").Value = 'Ilse Aguirre_20251120_155542'
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = 'Ilse Aguirre'
$ws.Range("D10").Value = 24
$ws.Range("E10").Value = 'Female'
$ws.Range("F10").Value = '2025-11-20 15:55:42'
$ws.Range("G10").Value = '{
  "portion": 0.4,
  "diet": 0.7142857142857143,
  "salt": 0.2,
  "fat": 0.6,
  "natural": 0.4,
  "convenience": 0.2,
  "price": 0.2
}'
$ws.Range("H10").Value = 'Nongshim Neoguri Spicy Seafood'
$c = $ws.Range("I10")
$c.NumberFormat = "@"
$c.Value = '0.533'
$c.Style = "Normal"
$ws.Range("J10").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K10").Value = 'Nissin Chow Mein Teriyaki Beef'
$c = $ws.Range("L10")
$c.NumberFormat = "@"
$c.Value = '0.422'
$c.Style = "Normal"
$ws.Range("M10").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N10").Value = 'Nongshim Shin Ramyun'
$c = $ws.Range("O10")
$c.NumberFormat = "@"
$c.Value = '0.419'
$c.Style = "Normal"
$ws.Range("P10").Value = 'Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio'
$ws.Range("Q10").Value = 'Amy’s Macaroni & Cheese (frozen)'
$c = $ws.Range("R10")
$c.NumberFormat = "@"
$c.Value = '0.631'
$c.Style = "Normal"
$ws.Range("S10").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("T10").Value = 'Kraft Macaroni & Cheese Dinner'
$c = $ws.Range("U10")
$c.NumberFormat = "@"
$c.Value = '0.622'
$c.Style = "Normal"
$ws.Range("V10").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("W10").Value = 'Annie’s Shells & White Cheddar'
$c = $ws.Range("X10")
$c.NumberFormat = "@"
$c.Value = '0.587'
$c.Style = "Normal"
$ws.Range("Y10").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("Z10").Value = 'Wild Planet Wild Tuna Pasta Salad'
$c = $ws.Range("AA10")
$c.NumberFormat = "@"
$c.Value = '0.762'
$c.Style = "Normal"
$ws.Range("AB10").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC10").Value = 'StarKist Chicken Creations (Chicken Salad)'
$c = $ws.Range("AD10")
$c.NumberFormat = "@"
$c.Value = '0.544'
$c.Style = "Normal"
$ws.Range("AE10").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AF10").Value = 'Kitchens of India Variety Pack'
$c = $ws.Range("AG10")
$c.NumberFormat = "@"
$c.Value = '0.497'
$c.Style = "Normal"
$ws.Range("AH10").Value = 'Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad'
